# Auto-generated edit script: updates the cryptos price table (rows 2-51)
# Source: GitHub Actions crypto price refresh, Sun Aug  6 08:25:07 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to store a literal text value (matches the workbook's
    # original inline-string cells) instead of letting Excel's input parser
    # reinterpret numeric-looking strings (e.g. "1.000", "23.29") as numbers.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2: D2, E2
Set-TextValue $ws.Range("D2") '29.129.00'
Set-TextValue $ws.Range("E2") '  +0.37%  '

# Row 3: D3, E3
Set-TextValue $ws.Range("D3") '1.836.80'
Set-TextValue $ws.Range("E3") '  +0.45%  '

# Row 4: D4, E4
Set-TextValue $ws.Range("D4") '1.000'
Set-TextValue $ws.Range("E4") '  +0.14%  '

# Row 5: D5, E5
Set-TextValue $ws.Range("D5") '245.44'
Set-TextValue $ws.Range("E5") '  +1.85%  '

# Row 6: D6, E6
Set-TextValue $ws.Range("D6") '0.6360'
Set-TextValue $ws.Range("E6") '  +2.07%  '

# Row 7: E7
Set-TextValue $ws.Range("E7") '  +0.17%  '

# Row 8: D8, E8
Set-TextValue $ws.Range("D8") '0.07560'
Set-TextValue $ws.Range("E8") '  +0.27%  '

# Row 9: D9, E9
Set-TextValue $ws.Range("D9") '0.2954'
Set-TextValue $ws.Range("E9") '  +1.63%  '

# Row 10: D10, E10
Set-TextValue $ws.Range("D10") '23.29'
Set-TextValue $ws.Range("E10") '  +2.53%  '

# Row 11: D11, E11
Set-TextValue $ws.Range("D11") '0.07728'
Set-TextValue $ws.Range("E11") '  +1.15%  '

# Row 12: D12, E12
Set-TextValue $ws.Range("D12") '1.835.90'
Set-TextValue $ws.Range("E12") '  +0.25%  '

# Row 13: D13, E13
Set-TextValue $ws.Range("D13") '5.012'
Set-TextValue $ws.Range("E13") '  +1.25%  '

# Row 14: D14, E14
Set-TextValue $ws.Range("D14") '0.6736'
Set-TextValue $ws.Range("E14") '  +1.51%  '

# Row 15: D15, E15
Set-TextValue $ws.Range("D15") '83.42'
Set-TextValue $ws.Range("E15") '  +1.48%  '

# Row 16: D16, E16
Set-TextValue $ws.Range("D16") '0.000009589'
Set-TextValue $ws.Range("E16") '  +5.61%  '

# Row 17: D17, E17
Set-TextValue $ws.Range("D17") '6.101'
Set-TextValue $ws.Range("E17") '  +2.16%  '

# Row 18: D18, E18
Set-TextValue $ws.Range("D18") '29.152.15'
Set-TextValue $ws.Range("E18") '  +0.72%  '

# Row 19: D19, E19
Set-TextValue $ws.Range("D19") '12.64'
Set-TextValue $ws.Range("E19") '  +2.67%  '

# Row 20: D20, E20
Set-TextValue $ws.Range("D20") '228.18'
Set-TextValue $ws.Range("E20") '  +1.72%  '

# Row 21: E21
Set-TextValue $ws.Range("E21") '  +0.08%  '

# Row 22: D22, E22
Set-TextValue $ws.Range("D22") '7.199'
Set-TextValue $ws.Range("E22") '  +0.17%  '

# Row 23: E23
Set-TextValue $ws.Range("E23") '  +0.10%  '

# Row 24: D24, E24
Set-TextValue $ws.Range("D24") '160.65'
Set-TextValue $ws.Range("E24") '  +0.58%  '

# Row 25: D25, E25
Set-TextValue $ws.Range("D25") '0.1421'
Set-TextValue $ws.Range("E25") '  +4.85%  '

# Row 26: D26, E26
Set-TextValue $ws.Range("D26") '8.568'
Set-TextValue $ws.Range("E26") '  +2.10%  '

# Row 27: D27, E27
Set-TextValue $ws.Range("D27") '17.98'
Set-TextValue $ws.Range("E27") '  +1.07%  '

# Row 28: D28, E28
Set-TextValue $ws.Range("D28") '1.506'
Set-TextValue $ws.Range("E28") '  +0.97%  '

# Row 29: D29, E29
Set-TextValue $ws.Range("D29") '4.168'
Set-TextValue $ws.Range("E29") '  +3.21%  '

# Row 30: D30, E30
Set-TextValue $ws.Range("D30") '4.078'
Set-TextValue $ws.Range("E30") '  +1.44%  '

# Row 31: B31, C31, D31, E31
Set-TextValue $ws.Range("B31") 'Hedera'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D31") '0.05437'
Set-TextValue $ws.Range("E31") '  +4.72%  '

# Row 32: B32, C32, D32, E32
Set-TextValue $ws.Range("B32") 'Toncoin'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D32") '1.200'
Set-TextValue $ws.Range("E32") '  -0.09%  '

# Row 33: E33
Set-TextValue $ws.Range("E33") '  +1.16%  '

# Row 34: D34, E34
Set-TextValue $ws.Range("D34") '0.7484'
Set-TextValue $ws.Range("E34") '  +2.75%  '

# Row 35: D35, E35
Set-TextValue $ws.Range("D35") '1.144'
Set-TextValue $ws.Range("E35") '  -0.52%  '

# Row 36: D36, E36
Set-TextValue $ws.Range("D36") '2.666'
Set-TextValue $ws.Range("E36") '  +1.97%  '

# Row 37: D37, E37
Set-TextValue $ws.Range("D37") '1.249.17'
Set-TextValue $ws.Range("E37") '  -2.27%  '

# Row 38: E38
Set-TextValue $ws.Range("E38") '  +0.49%  '

# Row 39: D39, E39
Set-TextValue $ws.Range("D39") '2.760'
Set-TextValue $ws.Range("E39") '  +0.05%  '

# Row 40: D40, E40
Set-TextValue $ws.Range("D40") '6.680'
Set-TextValue $ws.Range("E40") '  +4.55%  '

# Row 41: D41, E41
Set-TextValue $ws.Range("D41") '0.9060'
Set-TextValue $ws.Range("E41") '  +1.48%  '

# Row 42: D42, E42
Set-TextValue $ws.Range("D42") '1.003'
Set-TextValue $ws.Range("E42") '  +0.22%  '

# Row 43: D43, E43
Set-TextValue $ws.Range("D43") '101.68'
Set-TextValue $ws.Range("E43") '  +0.30%  '

# Row 44: D44, E44
Set-TextValue $ws.Range("D44") '1.986.73'
Set-TextValue $ws.Range("E44") '  +0.34%  '

# Row 45: D45, E45
Set-TextValue $ws.Range("D45") '0.00000000125'
Set-TextValue $ws.Range("E45") '  +4.59%  '

# Row 46: D46, E46
Set-TextValue $ws.Range("D46") '65.48'
Set-TextValue $ws.Range("E46") '  +3.28%  '

# Row 47: D47, E47
Set-TextValue $ws.Range("D47") '0.5115'
Set-TextValue $ws.Range("E47") '  +0.16%  '

# Row 48: D48, E48
Set-TextValue $ws.Range("D48") '0.4083'
Set-TextValue $ws.Range("E48") '  +2.81%  '

# Row 49: D49, E49
Set-TextValue $ws.Range("D49") '9.022'
Set-TextValue $ws.Range("E49") '  +1.91%  '

# Row 50: B50, C50, D50, E50
Set-TextValue $ws.Range("B50") 'Aptos'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D50") '6.784'
Set-TextValue $ws.Range("E50") '  +1.54%  '

# Row 51: B51, C51, D51, E51
Set-TextValue $ws.Range("B51") 'RenderToken'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D51") '1.652'
Set-TextValue $ws.Range("E51") '  +0.90%  '
